# Rerun all experiment using PR curve method to find threshold
$wb = $excel.ActiveWorkbook

# --- Sheet: "Top 10 players goal 90" ---
$ws = $wb.Worksheets.Item("Top 10 players goal 90")
$ws.Activate() | Out-Null
$ws.Range("D2").Value = 97
$ws.Range("D3").Value = 150
$ws.Range("D4").Value = 92
$ws.Range("D5").Value = 117
$ws.Range("D6").Value = 155
$ws.Range("D7").Value = 90
$ws.Range("D8").Value = 77
$ws.Range("D10").Value = 134
$ws.Range("D12").Select() | Out-Null

# --- Sheet: "Top 10 players assist 90" ---
$ws = $wb.Worksheets.Item("Top 10 players assist 90")
$ws.Activate() | Out-Null
$ws.Range("D2").Value = 90
$ws.Range("D3").Value = 5
$ws.Range("D5").Value = 150
$ws.Range("D6").Value = 84
$ws.Range("D7").Value = 93
$ws.Range("D8").Value = 78
$ws.Range("D9").Value = 95
$ws.Range("D10").Value = 77
$ws.Range("D11").Value = 117
$ws.Range("D10").Select() | Out-Null

# --- Sheet: "Top 10 players goal assist 90" ---
$ws = $wb.Worksheets.Item("Top 10 players goal assist 90")
$ws.Activate() | Out-Null
$ws.Range("D2").Value = 90
$ws.Range("D3").Value = 150
$ws.Range("D4").Value = 97
$ws.Range("D5").Value = 92
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 117
$ws.Range("D9").Value = 155
$ws.Range("D10").Value = 77
$ws.Range("D11").Value = 84
$ws.Range("D12").Select() | Out-Null

# --- Sheet: "Top 10 proposed ranking" (swap rows 6 and 7 in columns B and C) ---
$ws = $wb.Worksheets.Item("Top 10 proposed ranking")
$ws.Activate() | Out-Null
$ws.Range("B6").Value = "Jordi Alba"
$ws.Range("C6").Value = "9 million euro"
$ws.Range("B7").Value = "Jorge Resurreccion Merodio "
$ws.Range("C7").Value = "35 million euro"
$ws.Range("F7").Select() | Out-Null
